$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291140689345"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911428408551"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911428408551"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911429102745"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291142989275"

# Sheet1 (GNG) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911406334639.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911406584082.csv"
$ws1.Range("B4").Value = "go_stims-16502911406594048.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291140689345.csv"

# Sheet2 (NB) updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16502911426256652.csv"
$ws2.Range("B3").Value = "ZB-match_6-16502911410207195.csv"
$ws2.Range("B4").Value = "OB-16502911413472333.csv"
$ws2.Range("B5").Value = "OB-16502911424437702.csv"
$ws2.Range("B6").Value = "ZB-match_2-16502911411713958.csv"
$ws2.Range("B7").Value = "ZB-match_6-16502911412851353.csv"
$ws2.Range("B8").Value = "TB-16502911428203955.csv"
$ws2.Range("B9").Value = "OB-16502911420623624.csv"
$ws2.Range("B10").Value = "TB-16502911427230709.csv"

# Sheet3 (RS) - no cell changes

# Sheet4 (TOL) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911428643274.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911428442965.csv"
$ws4.Range("B4").Value = "MM_stims-1650291142895082.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291142865334.csv"
$ws4.Range("B6").Value = "MM_stims-16502911429102745.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911428960829.csv"

# Sheet5 (vSAT) updates
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650291142973375.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911429427254.csv"
$ws5.Range("B4").Value = "SAT_stims-16502911429140556.csv"
$ws5.Range("B5").Value = "SAT_stims-16502911429283297.csv"
